$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 561.6908063895213
$ws.Range("D2").Value = 133.3636251146518
$ws.Range("F2").Value = 454
$ws.Range("G2").Value = 532
$ws.Range("H2").Value = 635
$ws.Range("C3").Value = 41.00357317009417
$ws.Range("D3").Value = 5.038087600419884
$ws.Range("F3").Value = 38.05
$ws.Range("G3").Value = 40.14
$ws.Range("H3").Value = 44.18
$ws.Range("C4").Value = 1.418916017527934
$ws.Range("D4").Value = 2.246730416347243
$ws.Range("F4").Value = 0.53
$ws.Range("G4").Value = 1.06
$ws.Range("H4").Value = 1.87
$ws.Range("C5").Value = 323.0147241417824
$ws.Range("D5").Value = 10.23774700587594
$ws.Range("F5").Value = 316.93
$ws.Range("G5").Value = 324.96
$ws.Range("H5").Value = 331.39
$ws.Range("C6").Value = 21.38153595175066
$ws.Range("D6").Value = 2.018052522979711
$ws.Range("F6").Value = 19.78
$ws.Range("G6").Value = 21.15
$ws.Range("H6").Value = 22.38
$ws.Range("C7").Value = -76.36750440812331
$ws.Range("D7").Value = 22.84409436797995
$ws.Range("G7").Value = -75
$ws.Range("H7").Value = -56
$ws.Range("C8").Value = 7.705910763377934
$ws.Range("D8").Value = 6.887995004726862
$ws.Range("C9").Value = 9.321625468734405
$ws.Range("D9").Value = 1.689066655511287
$ws.Range("C10").Value = 867.8305051638014
$ws.Range("D10").Value = 0.4610278333245824
$ws.Range("C11").Value = 0.5569541734257892
$ws.Range("D11").Value = 0.5909813791365603
$ws.Range("C12").Value = 22.70173996112297
$ws.Range("D12").Value = 12.28042894373855
$ws.Range("C13").Value = 0.6727105086807948
$ws.Range("D13").Value = 0.748781494070802
$ws.Range("C14").Value = 1.826308060093058
$ws.Range("D14").Value = 1.665016986672359
$ws.Range("C15").Value = 93.76750440812313
$ws.Range("D15").Value = 22.84409436798105
$ws.Range("F15").Value = 73.40000000000001
$ws.Range("G15").Value = 92.40000000000001
$ws.Range("C16").Value = -85.64339072166538
$ws.Range("D16").Value = 20.62162584814459
$ws.Range("F16").Value = -101.6389203414338
$ws.Range("G16").Value = -85.79706163635328
$ws.Range("H16").Value = -66.7376019773414
$ws.Range("C17").Value = -77.93747995828738
$ws.Range("D17").Value = 25.21584035225402
$ws.Range("F17").Value = -92.46183611348224
$ws.Range("G17").Value = -75.41392685158225
$ws.Range("H17").Value = -56.2376019773414